# Update cryptocurrency price/volume data in the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.446.65"
$ws.Range("E2").Value = "  +1.52%  "
$ws.Range("D3").Value = "1.879.34"
$ws.Range("E3").Value = "  +2.03%  "
$ws.Range("E4").Value = "  +0.49%  "
$ws.Range("D5").Value = "'240.38"
$ws.Range("E5").Value = "  +3.74%  "
$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  +1.12%  "
$ws.Range("E7").Value = "  +0.47%  "
$ws.Range("E8").Value = "  +7.40%  "
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("E10").Value = "  +2.02%  "
$ws.Range("D11").Value = "'0.0991"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Value = "2.152.43"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").Value = "'11.64"
$ws.Range("E13").Value = "  +0.27%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").Value = "'0.688"
$ws.Range("E14").Value = "  +1.86%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.863.39"
$ws.Range("E15").Value = "  +2.17%  "
$ws.Range("E16").Value = "  +2.47%  "
$ws.Range("D17").Value = "35.445.20"
$ws.Range("E17").Value = "  +1.49%  "
$ws.Range("D18").Value = "'71.07"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("D19").Value = "0.0₃0805"
$ws.Range("E19").Value = "  +2.18%  "
$ws.Range("D20").Value = "'243.09"
$ws.Range("E20").Value = "  +1.32%  "
$ws.Range("D21").Value = "'12.38"
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").Value = "'4.77"
$ws.Range("E22").Value = "  +1.96%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  -0.66%  "
$ws.Range("D25").Value = "'170.37"
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "'8.27"
$ws.Range("E26").Value = "  +5.81%  "
$ws.Range("B27").Value = "PancakeSwap"
$ws.Range("C27").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D27").Value = "'1.90"
$ws.Range("E27").Value = "  +24.78%  "
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  +1.70%  "
$ws.Range("D30").Value = "'0.0564"
$ws.Range("E30").Value = "  +2.34%  "
$ws.Range("E31").Value = "  +2.81%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").Value = "'1.83"
$ws.Range("E33").Value = "  +24.51%  "
$ws.Range("E34").Value = "  +2.41%  "
$ws.Range("E35").Value = "  +7.11%  "
$ws.Range("E36").Value = "  +18.15%  "
$ws.Range("D37").Value = "'1.32"
$ws.Range("E37").Value = "  +7.68%  "
$ws.Range("E38").Value = "  +2.86%  "
$ws.Range("D39").Value = "'0.0205"
$ws.Range("E39").Value = "  +5.36%  "
$ws.Range("D40").Value = "'91.31"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("D41").Value = "1.355.59"
$ws.Range("E41").Value = "  +0.49%  "
$ws.Range("D42").Value = "'15.27"
$ws.Range("E42").Value = "  +3.17%  "
$ws.Range("D43").Value = "'0.0605"
$ws.Range("E43").Value = "  +15.66%  "
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").Value = "'13.09"
$ws.Range("E45").Value = "  +57.70%  "
$ws.Range("E46").Value = "  +0.56%  "
$ws.Range("E47").Value = "  +6.71%  "
$ws.Range("E48").Value = "  -1.27%  "
$ws.Range("D49").Value = "2.064.16"
$ws.Range("E49").Value = "  +2.03%  "
$ws.Range("D50").Value = "'0.0689"
$ws.Range("E50").Value = "  +3.40%  "
$ws.Range("E51").Value = "  +0.29%  "
